# #327 Ajout des profils d'acces a58d18c1e8091c98efec92c8c093b361a253eee5
#
# 1. Metadata sheet: bump the "Date" value.
# 2. Elements sheet: swap the two "Mapping" columns (AK <-> AL) - both the
#    header text and the per-row values - for rows 1, 3, 5 and 6.

$wb = $excel.ActiveWorkbook

$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsMetadata.Range("B8").Value = "2024-03-19T13:17:15+00:00"

$wsElements = $wb.Worksheets.Item("Elements")

function Swap-CellValues($ws, $addr1, $addr2) {
    $r1 = $ws.Range($addr1)
    $r2 = $ws.Range($addr2)
    $v1 = $r1.Value2()
    $v2 = $r2.Value2()
    $r1.Value = $v2
    $r2.Value = $v1
}

Swap-CellValues $wsElements "AK1" "AL1"
Swap-CellValues $wsElements "AK3" "AL3"
Swap-CellValues $wsElements "AK5" "AL5"
Swap-CellValues $wsElements "AK6" "AL6"
